# Updates the EC (Estado de Cuenta) database: removes the worker row for
# document 1047468017 (RAFAEL DE ZUBIRIA CABRALES) and refreshes the
# summary totals (Valor Mora, Cant. Trabajadores, Cant. Periodos) to
# reflect the remaining single worker (1143466687 - CAMILA SAUMETH PALOMINO).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the entire row for the worker that is being taken out of the
# statement; this shifts the remaining data (and the signature block
# below) up by one row, exactly like deleting a row in the Excel UI.
$ws.Rows.Item(16).Delete()

# Refresh the summary figures at the top of the statement.
$ws.Range("E11").Value = 58666
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1

# Column D ("Nombre Trabajador") was best-fit to the longest remaining
# name after the row removal; re-fit it to the new content width.
$ws.Columns.Item(4).ColumnWidth = 26.1666666666667
